$d = $word.ActiveDocument

$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"
$cr = [char]13

# ---------------------------------------------------------------------------
# 1. Insert a new "Meta description" paragraph right after the title
#    (Heading1) paragraph at the top of the document.
# ---------------------------------------------------------------------------
$titleText = "Play Celtic Charm Fire Blaze Quattro for Free - Review"

$titleParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $txt = $d.Paragraphs.Item($i).Range.Text.TrimEnd($cr)
    if ($txt -eq $titleText) {
        $titleParaIndex = $i
        break
    }
}

$titlePara = $d.Paragraphs.Item($titleParaIndex)
$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs.Item($titleParaIndex + 1)
$metaPara.Style = "Normal"

$metaXml = "<w:p xmlns:w='$wNs'><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t xml:space=`"preserve`">: Discover the exciting gameplay mechanics and features of Celtic Charm Fire Blaze Quattro. Play for free with multiple bonus rounds and free spins.</w:t></w:r></w:p>"
$null = $metaPara.Range.InsertXML($metaXml)

# ---------------------------------------------------------------------------
# 2. Remove the duplicated bold "Play Celtic Charm Fire Blaze Quattro for
#    Free - Review" paragraph that used to sit right before the closing
#    italic meta-description paragraph at the end of the document.
# ---------------------------------------------------------------------------
$dupParaIndex = -1
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $txt = $d.Paragraphs.Item($i).Range.Text.TrimEnd($cr)
    if ($txt -eq $titleText) {
        $dupParaIndex = $i
        break
    }
}

$dupTitlePara = $d.Paragraphs.Item($dupParaIndex)
$nextPara = $d.Paragraphs.Item($dupParaIndex + 1)
$dupRange = $d.Range($dupTitlePara.Range.Start, $nextPara.Range.Start)
$dupRange.Delete()

# ---------------------------------------------------------------------------
# 3. Replace the trailing italic paragraph's text with the new DALLE prompt.
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastRange = $lastPara.Range
$textLen = $lastRange.Text.TrimEnd($cr).Length
$oldTextRange = $d.Range($lastRange.Start, $lastRange.Start + $textLen)
$oldTextRange.Delete()

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$promptXml = "<w:p xmlns:w='$wNs'><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Prompt for DALLE: Create a colorful and engaging cartoon-style image for the game &quot;Celtic Charm Fire Blaze Quattro&quot; that features a happy Maya warrior with glasses. The image should convey a sense of excitement and good fortune, with the warrior surrounded by symbols of Celtic charm and luck. The warrior should be depicted as confident and joyful, with a big smile on their face and a twinkle in their eye. The image should also include the game's title and some of its key symbols, such as the four-leaf clover and the female wild card. Use bright colors and bold lines to make the image stand out and capture the attention of potential players.</w:t></w:r></w:p>"
$null = $lastPara.Range.InsertXML($promptXml)
